$d = $word.ActiveDocument

$replacements = @(
    @{old="872×3=2616"; new="397×4=1588"},
    @{old="774×6=4644"; new="945×3=2835"},
    @{old="694×3=2082"; new="931×9=8379"},
    @{old="214×2=428";  new="148×3=444"},
    @{old="359×8=2872"; new="963×6=5778"},
    @{old="974×9=8766"; new="672×8=5376"},
    @{old="294×8=2352"; new="907×3=2721"},
    @{old="829×2=1658"; new="556×5=2780"},
    @{old="387×5=1935"; new="317×2=634"},
    @{old="404×5=2020"; new="557×6=3342"},
    @{old="337×7=2359"; new="792×4=3168"},
    @{old="272×8=2176"; new="649×3=1947"},
    @{old="507×8=4056"; new="498×8=3984"},
    @{old="359×3=1077"; new="308×5=1540"},
    @{old="289×2=578";  new="459×7=3213"},
    @{old="199×5=995";  new="749×7=5243"},
    @{old="132×4=528";  new="302×4=1208"},
    @{old="875×4=3500"; new="337×8=2696"},
    @{old="312×5=1560"; new="248×8=1984"},
    @{old="828×4=3312"; new="185×7=1295"},
    @{old="125×6=750";  new="338×4=1352"},
    @{old="838×6=5028"; new="342×7=2394"},
    @{old="635×6=3810"; new="356×6=2136"},
    @{old="755×8=6040"; new="382×7=2674"},
    @{old="656×6=3936"; new="829×6=4974"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
